# Add a new "min_per_month" column (H) to the Tasks sheet, cleaned from
# functions input (values are now plain numbers instead of formulas).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

$ws.Range("H1").Value = "min_per_month"

$values = @(1,1,1,1,0,0,1,1,1,1,1,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
